$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2) | Out-Null
}

# Title heading (appears twice: main Heading1 and bold footer line)
Replace-Text "Play Flaming Fox Free - Stunning Graphics and High Volatility" "Play Flaming Fox for Free - Stunning Graphics & Exciting Bonus Features"

# "What we like" bullet list
Replace-Text "Impeccable attention to detail and stunning graphics create the perfect atmosphere" "Impeccable attention to detail with stunning graphics"
Replace-Text "Bonus system with randomly triggered special features adds excitement to gameplay" "Authentic oriental music creates an immersive atmosphere"
Replace-Text "Uncommon feature of winning combinations from both left to right and right to left" "Innovative Bonus system with random special features"
Replace-Text "High volatility and above average RTP make for a thrilling slot experience" "Possibility to win extra Free Spins during the Bonus feature"

# "What we don't like" bullet list
Replace-Text "Limited to only 10 fixed paylines" "High volatility may not appeal to players seeking frequent wins"
Replace-Text "Special features triggered by the ninja fox may be difficult to land" "Limited maximum bet of €20 may not satisfy high rollers"

# Italic summary line at the end
Replace-Text "Experience the stunning attention to detail of Flaming Fox and its high volatility gameplay for free. Trigger special features for extra excitement." "Read our review of Flaming Fox, a visually stunning slot game with innovative bonus features. Play for free now!"
